$wb = $excel.ActiveWorkbook

# Rename the existing sheet and add the new "TiempoConversion" sheet after it.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Uso de CPU"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TiempoConversion"

# Text/labels (order chosen to match shared-string insertion order of the target file).
$ws2.Range("B2").Value = "Worker 1"
$ws2.Range("B1").Value = "Tiempos de conversión muestreo"
$ws2.Range("C2").Value = "Worker 2"
$ws2.Range("D2").Value = "Worker 3"
$ws2.Range("H1").Value = "min"
$ws2.Range("H2").Value = "seg/min"
$ws2.Range("H3").Value = "Workers"
$ws2.Range("H4").Value = "Peticiones/worker"

# Merge & center the title across B1:D1.
$ws2.Range("B1:D1").Merge()
$ws2.Range("B1:D1").HorizontalAlignment = -4108  # xlCenter

# Helper figures (G1:G3).
$ws2.Range("G1").Value = 10
$ws2.Range("G2").Value = 60
$ws2.Range("G3").Value = 3

# Sampled conversion times for the three workers.
$dataB = @(3.535, 3.202, 3.283, 3.394, 3.89, 3.516, 3.38, 3.538, 3.388)
$dataC = @(3.706, 4.751, 3.41, 3.63, 3.611, 3.933, 3.749, 3.588, 3.575)
$dataD = @(3.519, 3.477, 3.69, 4.81, 3.554, 3.471, 3.339, 3.473, 3.489)

for ($i = 0; $i -lt 9; $i++) {
    $row = 4 + $i
    $ws2.Cells.Item($row, 2).Value = $dataB[$i]
    $ws2.Cells.Item($row, 3).Value = $dataC[$i]
    $ws2.Cells.Item($row, 4).Value = $dataD[$i]
}

# Geometric-mean summary row.
$ws2.Range("B14").Formula = "=GEOMEAN(B4:B12)"
$ws2.Range("C14").Formula = "=GEOMEAN(C4:C12)"
$ws2.Range("D14").Formula = "=GEOMEAN(D4:D12)"
$ws2.Range("F14").Formula = "=GEOMEAN(B14:D14)"

# Requests-per-worker projections.
$ws2.Range("G4").Formula = "=G2*G1/F14"
$ws2.Range("G5").Formula = "=G3*G2*G1/F14"

# Match the saved selection/active cell on the new sheet.
$ws2.Range("H4").Select()

Write-Output "edit complete"
